$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure all target cells remain plain text (matches source data which used
# inline/shared strings, not numeric types) by forcing a text number format
# before assigning each value.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "76.446.51"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.57%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.044.41"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +4.46%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "201.84"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.79%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "624.48"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +4.68%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.22%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.207"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +4.79%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "3.047.99"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +4.53%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.439"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.15%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.57%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.27"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +7.40%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.605.37"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +4.40%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "29.33"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +3.20%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "76.419.11"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.64%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000193"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +1.97%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.042.57"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +3.60%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.55"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +3.86%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "9.05"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +2.07%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "376.40"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.92%  "
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.86%  "
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +1.46%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +3.24%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.209.84"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +4.61%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "4.41"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +3.94%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.999"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.01%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.88"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +1.64%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +3.41%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.997"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.03%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.32"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +7.22%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +1.68%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "502.41"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.12%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +6.90%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.02%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "20.83"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +2.87%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "162.40"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -1.68%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.387"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +4.80%  "
$ws.Range("B39").NumberFormat = "@"
$ws.Range("B39").Value = "Kaspa"
$ws.Range("C39").NumberFormat = "@"
$ws.Range("C39").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.117"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +4.13%  "
$ws.Range("B40").NumberFormat = "@"
$ws.Range("B40").Value = "WhiteBITCoin"
$ws.Range("C40").NumberFormat = "@"
$ws.Range("C40").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "20.04"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +2.29%  "
$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = "Aave"
$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "191.17"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +4.57%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -3.76%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.01%  "
$ws.Range("B44").NumberFormat = "@"
$ws.Range("B44").Value = "Mantle"
$ws.Range("C44").NumberFormat = "@"
$ws.Range("C44").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.802"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +22.14%  "
$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = "RenderToken"
$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.16"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +3.36%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +7.06%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "42.15"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +5.07%  "
$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = "dogwifhat"
$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.53"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +7.73%  "
$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = "Stacks"
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.67"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.37%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.611"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +6.65%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.92"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +5.38%  "
